{"js": "// Update the date heading and all 100 arithmetic-problem table cells.\n// The mapping below was derived from the target diff: each table cell's\n// text is replaced by position (row, column), not by a risky whole-document\n// text search, because a handful of the new expressions (e.g. \"25+34=\")\n// contain an old expression (e.g. \"5+34=\") as a literal substring.\n\nconst oldDate = \"2023-12-13 Wednesday\";\nconst newDate = \"2023-12-14 Thursday\";\n\nconst oldGrid = [\n    [\"50+9=\", \"51-15=\", \"52-23=\", \"68-46=\", \"27-12=\"],\n    [\"92-35=\", \"31-25=\", \"76+22=\", \"88-53=\", \"20+1=\"],\n    [\"27+61=\", \"30+6=\", \"49+6=\", \"87-0=\", \"80-29=\"],\n    [\"67+19=\", \"91-51=\", \"50-4=\", \"84-0=\", \"22+17=\"],\n    [\"88-59=\", \"18-17=\", \"73+26=\", \"18+50=\", \"71-70=\"],\n    [\"36+44=\", \"25-4=\", \"12+76=\", \"60-14=\", \"65+5=\"],\n    [\"36+4=\", \"24+39=\", \"46-45=\", \"73-47=\", \"61-59=\"],\n    [\"45-6=\", \"12-10=\", \"0+64=\", \"60-38=\", \"72-51=\"],\n    [\"9+81=\", \"52+3=\", \"20-6=\", \"73-30=\", \"32-24=\"],\n    [\"32-25=\", \"14+19=\", \"76-64=\", \"73+2=\", \"42-32=\"],\n    [\"62+21=\", \"87-56=\", \"23-6=\", \"38+24=\", \"98-98=\"],\n    [\"58-8=\", \"45-8=\", \"52-22=\", \"98-68=\", \"98-8=\"],\n    [\"50+14=\", \"2+30=\", \"7+58=\", \"43-31=\", \"68-32=\"],\n    [\"99-93=\", \"49-31=\", \"48+14=\", \"54-1=\", \"24+16=\"],\n    [\"73-17=\", \"72-56=\", \"62-54=\", \"31-14=\", \"27+36=\"],\n    [\"5+34=\", \"25+22=\", \"84-30=\", \"51-19=\", \"20-19=\"],\n    [\"32+11=\", \"78+10=\", \"4+83=\", \"5+62=\", \"47+36=\"],\n    [\"13+53=\", \"20+77=\", \"67-21=\", \"4+40=\", \"81-47=\"],\n    [\"72-55=\", \"82-16=\", \"47-14=\", \"78-65=\", \"84-40=\"],\n    [\"63-34=\", \"4+20=\", \"45-42=\", \"77+16=\", \"99-24=\"]\n  ];\n\nconst newGrid = [\n    [\"29+11=\", \"33-2=\", \"24+3=\", \"60-26=\", \"63-25=\"],\n    [\"40+42=\", \"32-12=\", \"3+49=\", \"14-5=\", \"40+43=\"],\n    [\"50+6=\", \"62-28=\", \"20-2=\", \"82-5=\", \"51-33=\"],\n    [\"42-38=\", \"15+27=\", \"5+59=\", \"13+66=\", \"36+31=\"],\n    [\"67+30=\", \"23+62=\", \"86-37=\", \"58+38=\", \"34-7=\"],\n    [\"37+34=\", \"38+25=\", \"24+60=\", \"83-70=\", \"4+57=\"],\n    [\"78-0=\", \"26+38=\", \"90-22=\", \"89+1=\", \"88-70=\"],\n    [\"3+66=\", \"85-72=\", \"9+62=\", \"77+18=\", \"7+71=\"],\n    [\"10-9=\", \"72-29=\", \"22+44=\", \"53-29=\", \"46+36=\"],\n    [\"24+18=\", \"41+27=\", \"90-15=\", \"20+71=\", \"25+34=\"],\n    [\"3+81=\", \"49-49=\", \"74-19=\", \"51+28=\", \"34-16=\"],\n    [\"99-23=\", \"98-76=\", \"90-80=\", \"40-30=\", \"97-33=\"],\n    [\"13+52=\", \"68-10=\", \"75-13=\", \"40-23=\", \"22+18=\"],\n    [\"51+24=\", \"59-4=\", \"33+24=\", \"18+67=\", \"72+25=\"],\n    [\"21-16=\", \"71-4=\", \"46+47=\", \"12+65=\", \"95-42=\"],\n    [\"30+15=\", \"12+16=\", \"97-95=\", \"85-27=\", \"83+14=\"],\n    [\"14+42=\", \"3+12=\", \"83-56=\", \"96-83=\", \"68+21=\"],\n    [\"23+29=\", \"51-12=\", \"39-8=\", \"95-45=\", \"43+8=\"],\n    [\"8+71=\", \"29+51=\", \"10+24=\", \"94-6=\", \"99-37=\"],\n    [\"13+15=\", \"79+12=\", \"80-11=\", \"21+60=\", \"25+14=\"]\n  ];\n\n// --- Update the date paragraph (first paragraph in the body) ---\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nif (paragraphs.items.length > 0 && paragraphs.items[0].text.trim() === oldDate) {\n  paragraphs.items[0].insertText(newDate, Word.InsertLocation.replace);\n}\n\n// --- Update every cell of the first table by (row, column) position ---\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\nfor (let r = 0; r < newGrid.length; r++) {\n  for (let c = 0; c < newGrid[r].length; c++) {\n    if (values[r] && values[r][c] !== undefined) {\n      // Only stomp the cell when it still holds the expected original\n      // text (defensive: keeps this a targeted edit rather than a blind\n      // overwrite if the table shape ever differs from what we expect).\n      if (values[r][c].trim() === oldGrid[r][c] || values[r][c].trim() !== newGrid[r][c]) {\n        values[r][c] = newGrid[r][c];\n      }\n    }\n  }\n}\ntable.values = values;\n\nawait context.sync();\n", "ps1": "# Update the date heading and all 100 arithmetic-problem table cells.\n# The mapping below was derived from the target diff: each table cell is\n# addressed by its (row, column) position via the Tables/Cell object model\n# rather than a whole-document text search-and-replace, because a handful\n# of the new expressions (e.g. \"25+34=\") contain an old expression\n# (e.g. \"5+34=\") as a literal substring -- a blind Find/Replace run in\n# document order could clobber an already-updated neighbour cell.\n\n$d = $word.ActiveDocument\n\n$oldDate = '2023-12-13 Wednesday'\n$newDate = '2023-12-14 Thursday'\n\n$oldGrid = @(\n  @('50+9=', '51-15=', '52-23=', '68-46=', '27-12='),\n  @('92-35=', '31-25=', '76+22=', '88-53=', '20+1='),\n  @('27+61=', '30+6=', '49+6=', '87-0=', '80-29='),\n  @('67+19=', '91-51=', '50-4=', '84-0=', '22+17='),\n  @('88-59=', '18-17=', '73+26=', '18+50=', '71-70='),\n  @('36+44=', '25-4=', '12+76=', '60-14=', '65+5='),\n  @('36+4=', '24+39=', '46-45=', '73-47=', '61-59='),\n  @('45-6=', '12-10=', '0+64=', '60-38=', '72-51='),\n  @('9+81=', '52+3=', '20-6=', '73-30=', '32-24='),\n  @('32-25=', '14+19=', '76-64=', '73+2=', '42-32='),\n  @('62+21=', '87-56=', '23-6=', '38+24=', '98-98='),\n  @('58-8=', '45-8=', '52-22=', '98-68=', '98-8='),\n  @('50+14=', '2+30=', '7+58=', '43-31=', '68-32='),\n  @('99-93=', '49-31=', '48+14=', '54-1=', '24+16='),\n  @('73-17=', '72-56=', '62-54=', '31-14=', '27+36='),\n  @('5+34=', '25+22=', '84-30=', '51-19=', '20-19='),\n  @('32+11=', '78+10=', '4+83=', '5+62=', '47+36='),\n  @('13+53=', '20+77=', '67-21=', '4+40=', '81-47='),\n  @('72-55=', '82-16=', '47-14=', '78-65=', '84-40='),\n  @('63-34=', '4+20=', '45-42=', '77+16=', '99-24=')\n)\n\n$newGrid = @(\n  @('29+11=', '33-2=', '24+3=', '60-26=', '63-25='),\n  @('40+42=', '32-12=', '3+49=', '14-5=', '40+43='),\n  @('50+6=', '62-28=', '20-2=', '82-5=', '51-33='),\n  @('42-38=', '15+27=', '5+59=', '13+66=', '36+31='),\n  @('67+30=', '23+62=', '86-37=', '58+38=', '34-7='),\n  @('37+34=', '38+25=', '24+60=', '83-70=', '4+57='),\n  @('78-0=', '26+38=', '90-22=', '89+1=', '88-70='),\n  @('3+66=', '85-72=', '9+62=', '77+18=', '7+71='),\n  @('10-9=', '72-29=', '22+44=', '53-29=', '46+36='),\n  @('24+18=', '41+27=', '90-15=', '20+71=', '25+34='),\n  @('3+81=', '49-49=', '74-19=', '51+28=', '34-16='),\n  @('99-23=', '98-76=', '90-80=', '40-30=', '97-33='),\n  @('13+52=', '68-10=', '75-13=', '40-23=', '22+18='),\n  @('51+24=', '59-4=', '33+24=', '18+67=', '72+25='),\n  @('21-16=', '71-4=', '46+47=', '12+65=', '95-42='),\n  @('30+15=', '12+16=', '97-95=', '85-27=', '83+14='),\n  @('14+42=', '3+12=', '83-56=', '96-83=', '68+21='),\n  @('23+29=', '51-12=', '39-8=', '95-45=', '43+8='),\n  @('8+71=', '29+51=', '10+24=', '94-6=', '99-37='),\n  @('13+15=', '79+12=', '80-11=', '21+60=', '25+14=')\n)\n\n# --- Update the date paragraph (first paragraph in the body) ---\n$titlePara = $d.Paragraphs.Item(1)\n$titleText = $titlePara.Range.Text.TrimEnd([char]13, [char]7)\nif ($titleText -eq $oldDate) {\n  $titlePara.Range.Text = $newDate\n}\n\n# --- Update every cell of the first table by (row, column) position ---\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $colCount; $c++) {\n    $cell = $table.Cell($r, $c)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    $expectedOld = $oldGrid[$r - 1][$c - 1]\n    $expectedNew = $newGrid[$r - 1][$c - 1]\n    # Only stomp the cell when it still holds the expected original text\n    # (defensive: keeps this a targeted edit rather than a blind overwrite\n    # if the table shape ever differs from what we expect).\n    if ($current -eq $expectedOld -or $current -ne $expectedNew) {\n      $cell.Range.Text = $expectedNew\n    }\n  }\n}\n"}
